$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''25.974.84'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.07%  '

# Row 3
$ws.Range("D3").Value = '''1.833.78'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.88%  '

# Row 4
$ws.Range("E4").Value = '  -0.23%  '

# Row 5
$ws.Range("D5").Value = '''278.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -7.23%  '

# Row 6
$ws.Range("D6").Value = '''0.9986'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.19%  '

# Row 7
$ws.Range("D7").Value = '''0.5109'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.83%  '

# Row 8
$ws.Range("D8").Value = '''0.3493'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.56%  '

# Row 9
$ws.Range("D9").Value = '''44.72'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.18%  '

# Row 10
$ws.Range("D10").Value = '''0.06819'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.55%  '

# Row 11
$ws.Range("D11").Value = '''19.96'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -7.32%  '

# Row 12
$ws.Range("D12").Value = '''0.8058'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -9.18%  '

# Row 13
$ws.Range("D13").Value = '''0.07804'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.10%  '

# Row 14
$ws.Range("D14").Value = '''1.829.89'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.38%  '

# Row 15
$ws.Range("D15").Value = '''5.073'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.07%  '

# Row 16
$ws.Range("D16").Value = '''88.25'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.52%  '

# Row 17
$ws.Range("D17").Value = '''0.9979'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.26%  '

# Row 18
$ws.Range("D18").Value = '''14.16'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.62%  '

# Row 19
$ws.Range("D19").Value = '''0.000008066'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.90%  '

# Row 20
$ws.Range("D20").Value = '''0.9984'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.19%  '

# Row 21
$ws.Range("D21").Value = '''26.007.78'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.08%  '

# Row 22
$ws.Range("D22").Value = '''4.754'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.26%  '

# Row 23
$ws.Range("D23").Value = '''10.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.69%  '

# Row 24
$ws.Range("D24").Value = '''6.206'
$ws.Range("D24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = '''2.376'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.96%  '

# Row 26
$ws.Range("D26").Value = '''142.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.46%  '

# Row 27
$ws.Range("D27").Value = '''1.667'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.20%  '

# Row 28
$ws.Range("D28").Value = '''17.17'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.40%  '

# Row 29
$ws.Range("D29").Value = '''109.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.62%  '

# Row 30
$ws.Range("D30").Value = '''4.359'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.19%  '

# Row 31
$ws.Range("D31").Value = '''4.286'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.34%  '

# Row 32
$ws.Range("D32").Value = '''0.08783'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.50%  '

# Row 33
$ws.Range("D33").Value = '''0.04862'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.10%  '

# Row 34
$ws.Range("D34").Value = '''1.162'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.86%  '

# Row 35
$ws.Range("D35").Value = '''0.7288'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -10.14%  '

# Row 36
$ws.Range("D36").Value = '''2.867'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.75%  '

# Row 37
$ws.Range("D37").Value = '''3.198'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.06%  '

# Row 38
$ws.Range("D38").Value = '''0.9978'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.26%  '

# Row 39
$ws.Range("D39").Value = '''2.396'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -10.52%  '

# Row 40
$ws.Range("D40").Value = '''0.01850'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.03%  '

# Row 41
$ws.Range("D41").Value = '''0.5138'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -15.19%  '

# Row 42
$ws.Range("D42").Value = '''0.9473'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -11.04%  '

# Row 43
$ws.Range("D43").Value = '''117.32'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.91%  '

# Row 44
$ws.Range("D44").Value = '''6.258'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.35%  '

# Row 45
$ws.Range("D45").Value = '''8.002'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.67%  '

# Row 46
$ws.Range("D46").Value = '''0.9981'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.22%  '

# Row 47
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '''0.4510'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -14.62%  '

# Row 48
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = '''0.1363'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.36%  '

# Row 49
$ws.Range("D49").Value = '''9.303'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.60%  '

# Row 50
$ws.Range("E50").Value = '  -3.06%  '

# Row 51
$ws.Range("D51").Value = '''0.05918'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.27%  '
